$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Wipe the whole old table (values + formats): the new layout only spans
# through row 13 (old data went through row 14) and shuffles every column,
# so start from a clean slate rather than trying to patch cells in place.
# ---------------------------------------------------------------------------
$ws.Range("A1:K14").Clear()

# ---------------------------------------------------------------------------
# Build a one-off named style that mirrors the existing "Arial 9 / General"
# cell format (fontId 1) but WITHOUT an explicit number-format application,
# matching the new cellXfs entry introduced by this edit. We apply it to the
# new numeric-header cells then drop the named style again so only the xf
# record (not the cellStyle/cellStyleXfs bookkeeping) survives.
# ---------------------------------------------------------------------------
$headerStyle = $wb.Styles.Add("HeaderHelper9")
$headerStyle.Font.Size = 9
$headerStyle.Font.Name = "Arial"

# ---------------------------------------------------------------------------
# Row 1: new header labels.
# Columns A:E get the plain default style (no explicit formatting), columns
# F:K get the Arial-9/General header style.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,1).Value = "idx"
$ws.Cells.Item(1,2).Value = "idx2"
$ws.Cells.Item(1,3).Value = "Name"
$ws.Cells.Item(1,4).Value = "Date Start"
$ws.Cells.Item(1,5).Value = "Date End"

$ws.Cells.Item(1,6).Value = "(m3/s)"
$ws.Cells.Item(1,7).Value = "(MW1)"
$ws.Cells.Item(1,8).Value = "(MW2)"
$ws.Cells.Item(1,9).Value = "(GWh) Winter"
$ws.Cells.Item(1,10).Value = "(GWh) Summer"
$ws.Cells.Item(1,11).Value = "(GWh) Year"

$ws.Range("F1:K1").Style = "HeaderHelper9"

# Drop the helper named style now that its xf has been stamped onto F1:K1.
$wb.Styles("HeaderHelper9").Delete()

# ---------------------------------------------------------------------------
# Data rows 2..13 (one per power plant, shifted up one row vs. the old
# sheet and with two new leading id columns).
# ---------------------------------------------------------------------------
$data = @(
  @(1, 303000, "Obermatt", 1905, 1963, 11, 24.22, 22.51, 19.71, 75.040000000000006, 94.75),
  @(2, 303900, "Wisserlen", 1905, 1997, 0.24, 1, 1, 1.05, 2.15, 3.2),
  @(3, 303700, "Unteraa", 1921, 1994, 32, 54, 54, 46.7, 35, 81.7),
  @(4, 303600, "Kaiserstuhl", 1933, $null, 10, 10.3, 9, 5.3, 21, 26.3),
  @(5, 302700, "Stalden (Kloster)", 1941, $null, 0.32, 1.27, 1.1499999999999999, 1.2, 3.5, 4.7),
  @(6, 304000, "Eichi", 1957, $null, 12, 2.4, 2.2000000000000002, 6.68, 7.72, 14.4),
  @(7, 303800, "Hugschwendi", 1960, $null, 2, 14, 14, 18, 19, 37),
  @(8, 303200, "Dallenwil", 1962, 1987, 14.7, 6.02, 5.45, 6.47, 20.03, 26.5),
  @(9, 303100, "Obermatt-Nebenzentrale", 1963, $null, 11, 1.49, 1.36, 0.76, 3.03, 3.79),
  @(10, 302900, "Arni", 1966, $null, 1, 0.7, 0.66, 0.11, 1.54, 1.65),
  @(11, 302800, "Engelberg", 1967, $null, 1.4, 0.77, 0.76, 0.28999999999999998, 1.45, 1.74),
  @(12, 303550, "Hackeren", 2003, $null, 0.22, 1.2, 1.2, 0.8, 2.4, 3.2)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,1).Font.Size = 9
    $ws.Cells.Item($r,1).NumberFormat = "0"

    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,2).Font.Size = 9
    $ws.Cells.Item($r,2).NumberFormat = "0"

    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,3).Font.Size = 9

    $ws.Cells.Item($r,4).Value = $row[3]
    $ws.Cells.Item($r,4).Font.Size = 9
    $ws.Cells.Item($r,4).NumberFormat = "0"

    if ($null -ne $row[4]) {
        $ws.Cells.Item($r,5).Value = $row[4]
        $ws.Cells.Item($r,5).Font.Size = 9
        $ws.Cells.Item($r,5).NumberFormat = "0"
    }

    for ($c = 6; $c -le 11; $c++) {
        $ws.Cells.Item($r,$c).Value = $row[$c - 1]
        $ws.Cells.Item($r,$c).Font.Size = 9
        $ws.Cells.Item($r,$c).NumberFormat = "0.00"
    }

    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Selection, matching the edited sheetView.
# ---------------------------------------------------------------------------
$ws.Range("A2:K2").Select() | Out-Null
